# step 15 - don't write sort if empty
# Fill in the previously-empty C12 cell with a placeholder (5 spaces) so that
# a value is always written instead of leaving the cell out of the sort range,
# matching the style already used by the rest of row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the neighboring cell (B12) onto C12 so the new
# cell reuses the existing cell style instead of creating a new one.
$ws.Cells.Item(12, 2).Copy()
$ws.Cells.Item(12, 3).PasteSpecial(-4122)

# Write the placeholder value into C12.
$ws.Cells.Item(12, 3).Value = "     "

# Update the active selection to the newly written cell.
$ws.Range("C12").Select()
